$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.404.19'
$ws.Range('E2').Value = '  -0.91%  '

# Row 3
$ws.Range('D3').Value = '3.286.02'
$ws.Range('E3').Value = '  -0.52%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.71'
$ws.Range('E5').Value = '  +3.83%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.31'
$ws.Range('E6').Value = '  -3.11%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.08%  '

# Row 8
$ws.Range('D8').Value = '3.281.36'
$ws.Range('E8').Value = '  -0.41%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.570'
$ws.Range('E9').Value = '  -2.82%  '

# Row 10
$ws.Range('E10').Value = '  -5.42%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.570'
$ws.Range('E11').Value = '  -2.68%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.40'
$ws.Range('E12').Value = '  -2.22%  '

# Row 13
$ws.Range('E13').Value = '  -3.33%  '

# Row 14
$ws.Range('B14').Value = 'BitcoinCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '627.13'
$ws.Range('E14').Value = '  +1.26%  '

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.810.15'
$ws.Range('E15').Value = '  -0.76%  '

# Row 16
$ws.Range('E16').Value = '  -3.04%  '

# Row 17
$ws.Range('D17').Value = '65.500.71'
$ws.Range('E17').Value = '  -0.83%  '

# Row 18
$ws.Range('E18').Value = '  -0.01%  '

# Row 19
$ws.Range('D19').Value = '3.288.62'
$ws.Range('E19').Value = '  -0.65%  '

# Row 20
$ws.Range('E20').Value = '  -2.55%  '

# Row 21
$ws.Range('E21').Value = '  -0.73%  '

# Row 22
$ws.Range('E22').Value = '  -2.39%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.94'
$ws.Range('E23').Value = '  -2.27%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '99.64'
$ws.Range('E24').Value = '  -2.27%  '

# Row 25
$ws.Range('E25').Value = '  -0.06%  '

# Row 26
$ws.Range('E26').Value = '  -0.01%  '

# Row 27
$ws.Range('E27').Value = '  -0.70%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.32'
$ws.Range('E28').Value = '  -3.00%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '30.55'
$ws.Range('E29').Value = '  +1.00%  '

# Row 30
$ws.Range('E30').Value = '  -3.66%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.45'
$ws.Range('E31').Value = '  -0.37%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '566.22'
$ws.Range('E32').Value = '  +1.15%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.67'
$ws.Range('E33').Value = '  -9.76%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.82'
$ws.Range('E34').Value = '  -2.16%  '

# Row 35
$ws.Range('D35').Value = '3.835.68'
$ws.Range('E35').Value = '  +0.30%  '

# Row 36
$ws.Range('E36').Value = '  -1.64%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.07%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.44'
$ws.Range('E38').Value = '  -3.29%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.125'
$ws.Range('E39').Value = '  -2.90%  '

# Row 40
$ws.Range('E40').Value = '  +5.48%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '32.32'
$ws.Range('E41').Value = '  -4.54%  '

# Row 42
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.11'
$ws.Range('E42').Value = '  -6.04%  '

# Row 43
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0677'
$ws.Range('E43').Value = '  -6.39%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.58'
$ws.Range('E44').Value = '  -5.22%  '

# Row 45
$ws.Range('E45').Value = '  -1.98%  '

# Row 46
$ws.Range('E46').Value = '  -4.12%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.01'
$ws.Range('E47').Value = '  -4.85%  '

# Row 48
$ws.Range('E48').Value = '  +0.27%  '

# Row 49
$ws.Range('E49').Value = '  -2.13%  '

# Row 50
$ws.Range('E50').Value = '  -2.88%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '129.09'
$ws.Range('E51').Value = '  +5.24%  '
